# Rename the "Vehículos" sheet to "Vehiculos" (accent removed) and make
# it the active sheet in the workbook (previously "Estaciones" was active).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Vehículos")
$ws.Name = "Vehiculos"

# Activating this sheet updates workbookView/@activeTab and moves the
# tabSelected flag from the old active sheet's sheetView to this one.
$ws.Activate()
